$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Scanner" to "Psychiatry"
$ws.Name = "Psychiatry"

# Remove the last data row (row 41) which duplicated row 40's data
$ws.Rows.Item(41).Delete()
